$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C (canDrinkAlcohol) and E (birthday) are formatted as Text ("@") and
# hold values that otherwise look like a boolean ("false") or a date
# ("01/01/1899"). Writing such literals straight into .Value lets Excel's
# automatic type detection turn them into a real boolean/date, so instead we
# copy the existing text cell (C9/E9, already stored as text) and paste only
# its value into the target cell - this keeps the literal as plain text -
# and then restore the "@" text number format used throughout those columns.
$xlPasteValues = -4163

# Row 9: update existing values -> new min test case (age -1, height -0.1, birthday 01/01/1899)
$ws.Range("B9").Value = -1
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial($xlPasteValues)
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "01/01/1899"
$ws.Range("F9").Value = -0.1

# Row 10 (new): max test case (age 999, height 99.9, birthday 01/01/2999)
$ws.Range("A10").Value = "Bob"
$ws.Range("B10").Value = 999
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial($xlPasteValues)
$ws.Range("C10").NumberFormat = "@"
$ws.Range("D10").Value = "AppleJuice"
$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial($xlPasteValues)
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "01/01/2999"
$ws.Range("F10").Value = 99.9

# Row 11 (new): re-add the original valid row that used to be row 9
$ws.Range("A11").Value = "Bob"
$ws.Range("B11").Value = 15
$ws.Range("C9").Copy()
$ws.Range("C11").PasteSpecial($xlPasteValues)
$ws.Range("C11").NumberFormat = "@"
$ws.Range("D11").Value = "AppleJuice"
$ws.Range("E9").Copy()
$ws.Range("E11").PasteSpecial($xlPasteValues)
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "01/01/2000"
$ws.Range("F11").Value = 1.7

$excel.CutCopyMode = $false
$ws.Range("J15").Select()
